$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "ok"
$ws.Range("A4").Value = "ok"
$ws.Range("A5").Value = "ok"
